# Update "想去人数" (want-to-go count) figures across sheets, as published
# by the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 86
$ws1.Range("F3").Value  = 147
$ws1.Range("F5").Value  = 3327
$ws1.Range("F6").Value  = 1056
$ws1.Range("F8").Value  = 2123
$ws1.Range("F9").Value  = 1110
$ws1.Range("F12").Value = 1673
$ws1.Range("F16").Value = 96
$ws1.Range("F18").Value = 1587
$ws1.Range("F19").Value = 634
$ws1.Range("F20").Value = 726
$ws1.Range("F21").Value = 607
$ws1.Range("F22").Value = 12255
$ws1.Range("F23").Value = 12298
$ws1.Range("F24").Value = 909
$ws1.Range("F27").Value = 37
$ws1.Range("F29").Value = 370
$ws1.Range("F30").Value = 1921

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value = 32

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 86
$ws4.Range("F4").Value  = 147
$ws4.Range("F6").Value  = 3327
$ws4.Range("F7").Value  = 1056
$ws4.Range("F9").Value  = 2123
$ws4.Range("F10").Value = 1110
$ws4.Range("F13").Value = 1673
$ws4.Range("F19").Value = 96
$ws4.Range("F22").Value = 1587
$ws4.Range("F23").Value = 634
$ws4.Range("F24").Value = 726
$ws4.Range("F25").Value = 607
$ws4.Range("F26").Value = 12255
$ws4.Range("F27").Value = 12298
$ws4.Range("F28").Value = 909
$ws4.Range("F31").Value = 37
$ws4.Range("F33").Value = 370
$ws4.Range("F34").Value = 1921
$ws4.Range("F40").Value = 32
